$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.01211971538532453
$ws.Range("D2").Value = 0.01707208025167617
$ws.Range("E2").Value = 0.1005247122900812
$ws.Range("F2").Value = 0.7921481765657319
$ws.Range("G2").Value = 0.6420672401678331
$ws.Range("H2").Value = 0.7124296235065373
$ws.Range("I2").Value = 0.5955098145560314
$ws.Range("K2").Value = 1.018700704066106
$ws.Range("M2").Value = 0.3622971501789465
$ws.Range("N2").Value = 1.122057150459391
$ws.Range("C3").Value = 0.01081830537506789
$ws.Range("D3").Value = 0.01667546020302169
$ws.Range("E3").Value = 0.09510043081046859
$ws.Range("F3").Value = 0.7789089341877826
$ws.Range("G3").Value = 0.6291456830562083
$ws.Range("H3").Value = 0.7117507354775512
$ws.Range("I3").Value = 0.5863354090676012
$ws.Range("K3").Value = 0.8960080471348704
$ws.Range("M3").Value = 0.3250381352683789
$ws.Range("N3").Value = 1.136711964987771
$ws.Range("C4").Value = 0.01001426740725719
$ws.Range("D4").Value = 0.01643156810214563
$ws.Range("E4").Value = 0.09185761988449315
$ws.Range("F4").Value = 0.7714106329806896
$ws.Range("G4").Value = 0.6217884110759968
$ws.Range("H4").Value = 0.7117919919531204
$ws.Range("I4").Value = 0.5811799721503377
$ws.Range("K4").Value = 0.8206849915816008
$ws.Range("M4").Value = 0.3022746997523811
$ws.Range("N4").Value = 1.146212808615154
$ws.Range("C5").Value = 0.009685380221391426
$ws.Range("D5").Value = 0.01633209899051025
$ws.Range("E5").Value = 0.09055794972849895
$ws.Range("F5").Value = 0.7685130055463034
$ws.Range("G5").Value = 0.6189345023970816
$ws.Range("H5").Value = 0.7119237628841688
$ws.Range("I5").Value = 0.5791987481297625
$ws.Range("K5").Value = 0.7899935067034107
$ws.Range("M5").Value = 0.2930267870296959
$ws.Range("N5").Value = 1.150210805081191
$ws.Range("C6").Value = 0.009630694561288067
$ws.Range("D6").Value = 0.01631557764664748
$ws.Range("E6").Value = 0.09034344978485009
$ws.Range("F6").Value = 0.76804138178926
$ws.Range("G6").Value = 0.6184693026840762
$ws.Range("H6").Value = 0.7119525797711788
$ws.Range("I6").Value = 0.5788769810643473
$ws.Range("K6").Value = 0.7848974238606843
$ws.Range("M6").Value = 0.2914928851329535
$ws.Range("N6").Value = 1.150882296588208
$ws.Range("C7").Value = 0.01000983690285295
$ws.Range("D7").Value = 0.0164302269407024
$ws.Range("E7").Value = 0.09184000414175131
$ws.Range("F7").Value = 0.771370915647708
$ws.Range("G7").Value = 0.6217493392475859
$ws.Range("H7").Value = 0.7117933039124011
$ws.Range("I7").Value = 0.5811527687885274
$ws.Range("K7").Value = 0.8202710614646946
$ws.Range("M7").Value = 0.3021498645907457
$ws.Range("N7").Value = 1.146266215627218
$ws.Range("C8").Value = 0.01167202797388711
$ws.Range("D8").Value = 0.01693540791134751
$ws.Range("E8").Value = 0.09863605748046211
$ws.Range("F8").Value = 0.7874520328505952
$ws.Range("G8").Value = 0.6374917496146963
$ws.Range("H8").Value = 0.7121003216904143
$ws.Range("I8").Value = 0.5922470544803176
$ws.Range("K8").Value = 0.9763943517736777
$ws.Range("M8").Value = 0.3494264810842225
$ws.Range("N8").Value = 1.127005735332283
$ws.Range("C9").Value = 0.01489171552979229
$ws.Range("D9").Value = 0.01792274394554738
$ws.Range("E9").Value = 0.1126714482458411
$ws.Range("F9").Value = 0.8240204585298301
$ws.Range("G9").Value = 0.6729744145181087
$ws.Range("H9").Value = 0.7163488456749718
$ws.Range("I9").Value = 0.6178159261746856
$ws.Range("K9").Value = 1.282627414910678
$ws.Range("M9").Value = 0.4430537523078755
$ws.Range("N9").Value = 1.093228640919374
$ws.Range("C10").Value = 0.01723253336157171
$ws.Range("D10").Value = 0.01864563538246955
$ws.Range("E10").Value = 0.1234332704963492
$ws.Range("F10").Value = 0.8540007717681277
$ws.Range("G10").Value = 0.7019094852435757
$ws.Range("H10").Value = 0.7217110258773403
$ws.Range("I10").Value = 0.6389604283844719
$ws.Range("K10").Value = 1.50767379723311
$ws.Range("M10").Value = 0.5124298405144145
$ws.Range("N10").Value = 1.070851341865232
$ws.Range("C11").Value = 0.01829200534367459
$ws.Range("D11").Value = 0.01897386116064581
$ws.Range("E11").Value = 0.1284308254929769
$ws.Range("F11").Value = 0.8683256343925905
$ws.Range("G11").Value = 0.7157070944033421
$ws.Range("H11").Value = 0.7246409625242336
$ws.Range("I11").Value = 0.6490995033004765
$ws.Range("K11").Value = 1.61006950256666
$ws.Range("M11").Value = 0.5441251498361197
$ws.Range("N11").Value = 1.061201992640001
$ws.Range("C12").Value = 0.01869241596086368
$ws.Range("D12").Value = 0.01909805345756155
$ws.Range("E12").Value = 0.1303382333213108
$ws.Range("F12").Value = 0.8738495566873894
$ws.Range("G12").Value = 0.7210240858241264
$ws.Range("H12").Value = 0.7258213061034837
$ws.Range("I12").Value = 0.6530142930409824
$ws.Range("K12").Value = 1.648846920119638
$ws.Range("M12").Value = 0.5561472483725396
$ws.Range("N12").Value = 1.057624407234115
$ws.Range("C13").Value = 0.0186062156754474
$ws.Range("D13").Value = 0.01907131101969739
$ws.Range("E13").Value = 0.1299267698433937
$ws.Range("F13").Value = 0.8726554493662775
$ws.Range("G13").Value = 0.7198748654174381
$ws.Range("H13").Value = 0.7255639427991412
$ws.Range("I13").Value = 0.6521678138798421
$ws.Range("K13").Value = 1.64049541321549
$ws.Range("M13").Value = 0.553557190953498
$ws.Range("N13").Value = 1.05839150292671
$ws.Range("C14").Value = 0.01832496325382493
$ws.Range("D14").Value = 0.01898408059126666
$ws.Range("E14").Value = 0.1285874479949456
$ws.Range("F14").Value = 0.8687780947017814
$ws.Range("G14").Value = 0.7161426751258375
$ws.Range("H14").Value = 0.7247366487198406
$ws.Range("I14").Value = 0.6494200624612034
$ws.Range("K14").Value = 1.613259699791286
$ws.Range("M14").Value = 0.5451138174054222
$ws.Range("N14").Value = 1.060906130081101
$ws.Range("C15").Value = 0.0181525849023032
$ws.Range("D15").Value = 0.01893063615109014
$ws.Range("E15").Value = 0.1277690290183884
$ws.Range("F15").Value = 0.8664160686413283
$ws.Range("G15").Value = 0.7138686246074997
$ws.Range("H15").Value = 0.7242391411784865
$ws.Range("I15").Value = 0.6477468125780774
$ws.Range("K15").Value = 1.596577331154435
$ws.Range("M15").Value = 0.5399445897463835
$ws.Range("N15").Value = 1.062456368367599
$ws.Range("C16").Value = 0.01716318486687385
$ws.Range("D16").Value = 0.01862417166515939
$ws.Range("E16").Value = 0.1231087440573546
$ws.Range("F16").Value = 0.8530784826557039
$ws.Range("G16").Value = 0.701020626448468
$ws.Range("H16").Value = 0.7215294451235934
$ws.Range("I16").Value = 0.6383083307368338
$ws.Range("K16").Value = 1.500982365048969
$ws.Range("M16").Value = 0.5103612308181766
$ws.Range("N16").Value = 1.071492632996211
$ws.Range("C17").Value = 0.01655483182216244
$ws.Range("D17").Value = 0.01843599908247739
$ws.Range("E17").Value = 0.1202761216803268
$ws.Range("F17").Value = 0.8450726504679125
$ws.Range("G17").Value = 0.6933020044620122
$ws.Range("H17").Value = 0.7199930044569243
$ws.Range("I17").Value = 0.6326517663075677
$ws.Range("K17").Value = 1.44234287029866
$ws.Range("M17").Value = 0.4922477495107387
$ws.Range("N17").Value = 1.07717200500672
$ws.Range("C18").Value = 0.01620441795366645
$ws.Range("D18").Value = 0.01832770932508154
$ws.Range("E18").Value = 0.1186564510547967
$ws.Range("F18").Value = 0.8405324999634729
$ws.Range("G18").Value = 0.6889221667496486
$ws.Range("H18").Value = 0.7191554499124635
$ws.Range("I18").Value = 0.6294472052726618
$ws.Range("K18").Value = 1.408616995314503
$ws.Range("M18").Value = 0.4818421192007349
$ws.Range("N18").Value = 1.080488523419188
$ws.Range("C19").Value = 0.0160856875173323
$ws.Range("D19").Value = 0.01829103465495052
$ws.Range("E19").Value = 0.1181096950196832
$ws.Range("F19").Value = 0.8390063565943535
$ws.Range("G19").Value = 0.6874494586313062
$ws.Range("H19").Value = 0.7188797886720124
$ws.Range("I19").Value = 0.6283705865973417
$ws.Range("K19").Value = 1.397198358392927
$ws.Range("M19").Value = 0.4783211364974278
$ws.Range("N19").Value = 1.081620003758239
$ws.Range("C20").Value = 0.01661964444728881
$ws.Range("D20").Value = 0.01845603643954874
$ws.Range("E20").Value = 0.1205766657572482
$ws.Range("F20").Value = 0.8459181950818504
$ws.Range("G20").Value = 0.6941174791745794
$ws.Range("H20").Value = 0.7201517811738682
$ws.Range("I20").Value = 0.6332488478500053
$ws.Range("K20").Value = 1.448584941052445
$ws.Range("M20").Value = 0.4941746376054681
$ws.Range("N20").Value = 1.076562261102577
$ws.Range("C21").Value = 0.01840759542514547
$ws.Range("D21").Value = 0.01900970506914135
$ws.Range("E21").Value = 0.1289804315009349
$ws.Range("F21").Value = 0.869914264537428
$ws.Range("G21").Value = 0.7172364026339721
$ws.Range("H21").Value = 0.7249777201820962
$ws.Range("I21").Value = 0.6502250947505956
$ws.Range("K21").Value = 1.6212594292316
$ws.Range("M21").Value = 0.5475933022608785
$ws.Range("N21").Value = 1.060165447877914
$ws.Range("C22").Value = 0.019571522251006
$ws.Range("D22").Value = 0.01937097388215037
$ws.Range("E22").Value = 0.1345600053021769
$ws.Range("F22").Value = 0.8861768184587646
$ws.Range("G22").Value = 0.7328833877236605
$ws.Range("H22").Value = 0.7285447947961359
$ws.Range("I22").Value = 0.661759447410077
$ws.Range("K22").Value = 1.734125813779826
$ws.Range("M22").Value = 0.5826208568593074
$ws.Range("N22").Value = 1.049894664923542
$ws.Range("C23").Value = 0.01895073868854524
$ws.Range("D23").Value = 0.01917821486014049
$ws.Range("E23").Value = 0.1315740066919702
$ws.Range("F23").Value = 0.8774439207829658
$ws.Range("G23").Value = 0.7244828443732558
$ws.Range("H23").Value = 0.7266030916571822
$ws.Range("I23").Value = 0.6555629697547118
$ws.Range("K23").Value = 1.673885858761196
$ws.Range("M23").Value = 0.5639153523312643
$ws.Range("N23").Value = 1.055335554800109
$ws.Range("C24").Value = 0.01659034472610443
$ws.Range("D24").Value = 0.01844697788203931
$ws.Range("E24").Value = 0.120440762263172
$ws.Range("F24").Value = 0.8455357296955555
$ws.Range("G24").Value = 0.6937486233830725
$ws.Range("H24").Value = 0.7200798556753085
$ws.Range("I24").Value = 0.632978759444029
$ws.Range("K24").Value = 1.445762942095371
$ws.Range("M24").Value = 0.4933034664454539
$ws.Range("N24").Value = 1.076837766302713
$ws.Range("C25").Value = 0.01402501266855438
$ws.Range("D25").Value = 0.01765604955195244
$ws.Range("E25").Value = 0.1087967959380975
$ws.Range("F25").Value = 0.8135842302260841
$ws.Range("G25").Value = 0.6628760666404503
$ws.Range("H25").Value = 0.7148072963177015
$ws.Range("I25").Value = 0.6104870984831479
$ws.Range("K25").Value = 1.199774077250026
$ws.Range("M25").Value = 0.4176238886613532
$ws.Range("N25").Value = 1.101938339040213
